$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.3763611233319027
$ws.Range("C4").Value = 0.322
$ws.Range("D4").Value = 0.6337611749680715
$ws.Range("E4").Value = 0.5620000000000001
$ws.Range("J4").Value = 0.6186990012492174
$ws.Range("K4").Value = 0.6129999999999999
$ws.Range("L4").Value = 0.6470572502149308
$ws.Range("M4").Value = 0.6260000000000001
$ws.Range("B5").Value = 0.6099383024474199
$ws.Range("C5").Value = 0.6040000000000001
$ws.Range("D5").Value = 0.6528329872542478
$ws.Range("E5").Value = 0.635
$ws.Range("F5").Value = 0.6584105147183463
$ws.Range("G5").Value = 0.952
$ws.Range("H5").Value = 0.503466578657865
$ws.Range("J5").Value = 0.6296802397477146
$ws.Range("K5").Value = 0.6449999999999999
$ws.Range("L5").Value = 0.6342670187606314
$ws.Range("M5").Value = 0.6239999999999999
$ws.Range("B6").Value = 0.3069578190985829
$ws.Range("C6").Value = 0.229
$ws.Range("D6").Value = 0.5945945787958171
$ws.Range("E6").Value = 0.579
$ws.Range("J6").Value = 0.6413186681572294
$ws.Range("K6").Value = 0.643
$ws.Range("L6").Value = 0.6568547824988145
$ws.Range("M6").Value = 0.6420000000000001
$ws.Range("B7").Value = 0.5735082173823869
$ws.Range("C7").Value = 0.58
$ws.Range("D7").Value = 0.6027026133312597
$ws.Range("E7").Value = 0.6110000000000001
$ws.Range("J7").Value = 0.5387609769255614
$ws.Range("K7").Value = 0.5629999999999999
$ws.Range("L7").Value = 0.5482004306528968
$ws.Range("M7").Value = 0.544
